# ------------------------------------------------------------------
# Edit: add "Player Info" and "ODI Batting Extra" sheets, rename the
# MATCH_CARD_LINK columns to MATCH_CODE (storing just the numeric
# match code instead of the full URL) on the existing "ODI Batting"
# and "ODI Bowling" sheets, and drop a few stray empty INNING_NUMBER
# cells.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper: pull the trailing numeric match code out of a howstat URL
function Get-MatchCode($url) {
    $parts = $url.Split("=")
    return $parts[$parts.Length - 1]
}

# ====================================================================
# 1) "ODI Batting" sheet: MATCH_CARD_LINK (col D) -> MATCH_CODE
# ====================================================================
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"
$battingLastRow = 89
$batting.Range("D2:D" + $battingLastRow).NumberFormat = "@"

for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -ne $null -and $url -ne "") {
        $cell.Value = Get-MatchCode $url
    }

    # drop the stray empty INNING_NUMBER cells (col B) entirely
    $inningCell = $batting.Cells.Item($r, 2)
    $inningVal = $inningCell.Value2
    if ($inningVal -eq $null -or $inningVal -eq "") {
        $inningCell.ClearContents()
    }
}

# ====================================================================
# 2) "ODI Bowling" sheet: MATCH_CARD_LINK (col B) -> MATCH_CODE
# ====================================================================
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"
$bowlingLastRow = 28
$bowling.Range("B2:B" + $bowlingLastRow).NumberFormat = "@"

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url -ne $null -and $url -ne "") {
        $cell.Value = Get-MatchCode $url
    }
}

# ====================================================================
# 3) New "ODI Batting Extra" sheet, appended after "ODI Bowling"
#    (done before the "Player Info" insertion below, since inserting
#    a sheet earlier in the collection can leave older sheet
#    references stale in this COM engine)
# ====================================================================
$battingExtra = $wb.Worksheets.Add($null, $bowling)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# force text formatting for the text-typed columns (A, C, D, E, F); B
# (BATTING_POSITION) is explicitly left as General (but still given a
# style) so numbers are stored as numbers, and blank cells still persist
# as empty placeholders instead of disappearing entirely
$battingExtra.Range("A2:A21").NumberFormat = "@"
$battingExtra.Range("B2:B21").NumberFormat = "General"
$battingExtra.Range("C2:F21").NumberFormat = "@"

$extraRows = @(
    @("4340", 1,   "5",  "0",  "15.42%", "NO"),
    @("4348", 2,   "10", "0",  "21.53%", "NO"),
    @("4377", 3,   "6",  "1",  "31.44%", "NO"),
    @("4378", 3,   "4",  "0",  "16.50%", "NO"),
    @("4379", 3,   "0",  "1",  "4.02%",  "NO"),
    @("4444", 3,   "0",  "0",  "1.05%",  "NO"),
    @("4446", 3,   "10", "2",  "39.62%", "YES"),
    @("4448", 3,   "0",  "0",  "0.75%",  "NO"),
    @("4525", $null, $null, $null, $null, "NO"),
    @("4528", 3,   "3",  "0",  "14.77%", "NO"),
    @("4530", 3,   "3",  "2",  "18.90%", "NO"),
    @("4537", 3,   "3",  "0",  "15.81%", "NO"),
    @("4538", $null, $null, $null, $null, "NO"),
    @("4539", 3,   "3",  "0",  "24.35%", "NO"),
    @("4582", 3,   "7",  "3",  "34.06%", "YES"),
    @("4585", 3,   "9",  "0",  "38.43%", "NO"),
    @("4588", 3,   "2",  "0",  "12.41%", "NO"),
    @("4671", $null, $null, $null, $null, "NO"),
    @("4674", 3,   "7",  "0",  "25.44%", "NO"),
    @("4675", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $battingExtra.Cells.Item($r, 1).Value = $row[0]

    if ($row[1] -ne $null) {
        $battingExtra.Cells.Item($r, 2).Value = $row[1]
    } else {
        $battingExtra.Cells.Item($r, 2).Value = ""
    }

    if ($row[2] -ne $null) {
        $battingExtra.Cells.Item($r, 3).Value = $row[2]
    } else {
        $battingExtra.Cells.Item($r, 3).Value = ""
    }

    if ($row[3] -ne $null) {
        $battingExtra.Cells.Item($r, 4).Value = $row[3]
    } else {
        $battingExtra.Cells.Item($r, 4).Value = ""
    }

    if ($row[4] -ne $null) {
        $battingExtra.Cells.Item($r, 5).Value = $row[4]
    } else {
        $battingExtra.Cells.Item($r, 5).Value = ""
    }

    $battingExtra.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}

# ====================================================================
# 4) New "Player Info" sheet, inserted before "ODI Batting"
#    (re-fetch the "ODI Batting" reference now that the sheet
#    collection has changed, to avoid a stale reference)
# ====================================================================
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingForInsert)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $playerInfoHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4009"
$playerInfo.Cells.Item(2, 2).Value = "Rahmat Shah"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Leg Break"

Write-Host "edit complete"
